# feat: add 2022-Q3 data
#
# - duplicate the existing "2022-Q1" sheet so its data survives unchanged
#   under the "2022-Q1" tab (now sitting after the new quarter's tab)
# - turn the original sheet into "2022-Q3" and fill it with the new
#   quarter's fund table
# - add a new row on the "总计" (summary) sheet for 2022-Q3, pushing the
#   existing 2022-Q1 row down

$xlPasteFormats = -4122
$xlCenter = -4108
$xlTop = -4160
$xlContinuous = 1

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$q1 = $wb.Worksheets.Item(2)

# 1) Duplicate "2022-Q1" right after itself; the copy becomes the new
#    "2022-Q1" tab, preserving its data & formatting untouched.
$q1.Copy($null, $q1)
$q1Copy = $wb.Worksheets.Item(3)

# 2) Wipe the original sheet (values + formatting) and rebuild it as the
#    2022-Q3 table.
$q1.Cells.Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "008928"
$q1.Range("C2").Value = "泰达宏利中证主要消费红利指数A"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "3.45"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "93.07"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "3.73"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.1287"
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "008929"
$q1.Range("C3").Value = "泰达宏利中证主要消费红利指数C"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "1.69"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "93.07"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "3.73"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.0630"
$q1.Range("H3").Value = 10

# Re-apply the bold/centered/bordered header & index-column look used
# elsewhere in this workbook (copied from the summary sheet's own header
# styling) to the freshly rebuilt table.
$summary.Range("B1:D1").Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)
$summary.Range("A2").Copy()
$q1.Range("A2:A3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# 3) Rename the sheets into their final positions/names.
$q1.Name = "2022-Q3"
$q1Copy.Name = "2022-Q1"

# 4) Update the summary sheet: shift the existing 2022-Q1 row down to row 3
#    and insert the new 2022-Q3 figures in row 2.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.52

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.19
